# Round the numeric data values (columns B:E, rows 2:13) to whole
# integers, matching the "write integer data" format change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("B2:E13")

for ($r = 1; $r -le $dataRange.Rows.Count; $r++) {
    for ($c = 1; $c -le $dataRange.Columns.Count; $c++) {
        $cell = $dataRange.Cells.Item($r, $c)
        $raw = [double]$cell.Value2
        $cell.Value = [Math]::Round($raw, 0)
    }
}
